# Add new Estonian rail line rows (40-43, 110) to the "Lines detail" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows, in the order they must appear (drives shared-string append order).
$rows = @(
    @{ Row = 55; A = 40;  B = "Türi—Paide—Tamsalu"; C = 1900; D = "750 mm";  E = "Dismantled";  F = 1972; G = "Extension to Tamsalu built 1915" },
    @{ Row = 56; A = 41;  B = "Narva—Musta";         C = 1969; D = "1520 mm"; E = "Dismantled";  F = 2001; G = "Rail line between Auvere and power plant dismantled" },
    @{ Row = 57; A = 42;  B = "Sonda—Mustvee";       C = 1926; D = "750 mm";  E = "Dismantled";  F = 1973; G = "~2km rebuilt as a museum railway" },
    @{ Row = 58; A = 43;  B = "Rakvere—Kunda";       C = 1870; D = "1520 mm"; E = "Freight only"; F = 2019; G = $null },
    @{ Row = 59; A = 110; B = "Tallinn—Narva";       C = 1870; D = "1520 mm"; E = "Operational";  F = 2021; G = $null }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the formatting (style) of the row above so the new row matches the
    # existing data rows (style index 1) without introducing new styles/fonts.
    $ws.Range("A" + ($rowNum - 1) + ":H" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum + ":H" + $rowNum).PasteSpecial(-4122)

    # Match the row height used throughout the rest of the sheet.
    $ws.Rows.Item($rowNum).RowHeight = 15.75

    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Value = $r.D
    $ws.Range("E" + $rowNum).Value = $r.E
    $ws.Range("F" + $rowNum).Value = $r.F
    if ($r.G -ne $null) {
        $ws.Range("G" + $rowNum).Value = $r.G
    } else {
        $ws.Range("G" + $rowNum).Clear()
    }
    $ws.Range("H" + $rowNum).Value = $false
}

$excel.CutCopyMode = $false

# Restore the view: scrolled so row 31 is at the top, with C59 selected.
$ws.Range("C59").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
